$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-parsed as a number by Excel
# are forced to the "@" (Text) number format first so the literal string
# (matching the original price-text formatting) is preserved verbatim.

$ws.Range("D2").Value = "57.429.20"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "2.317.21"
$ws.Range("E3").Value = "  -2.14%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.51"
$ws.Range("E5").Value = "  +2.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.43"
$ws.Range("E6").Value = "  -2.40%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -1.26%  "

$ws.Range("D9").Value = "2.339.61"
$ws.Range("E9").Value = "  -2.09%  "

$ws.Range("E10").Value = "  -1.48%  "

$ws.Range("E11").Value = "  +0.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.38"
$ws.Range("E12").Value = "  -2.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.352"
$ws.Range("E13").Value = "  +2.31%  "

$ws.Range("D14").Value = "2.734.95"
$ws.Range("E14").Value = "  -2.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.31"
$ws.Range("E15").Value = "  -4.40%  "

$ws.Range("D16").Value = "57.255.21"
$ws.Range("E16").Value = "  -0.34%  "

$ws.Range("E17").Value = "  -2.32%  "

$ws.Range("D18").Value = "2.328.72"
$ws.Range("E18").Value = "  -2.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "336.99"
$ws.Range("E19").Value = "  +2.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.39"
$ws.Range("E20").Value = "  -1.96%  "

$ws.Range("E22").Value = "  +0.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.81"
$ws.Range("E24").Value = "  +0.46%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.166"
$ws.Range("E25").Value = "  +0.27%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.70"
$ws.Range("E26").Value = "  -2.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.994"
$ws.Range("E27").Value = "  -0.27%  "

$ws.Range("E28").Value = "  +0.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.32"
$ws.Range("E29").Value = "  +3.76%  "

$ws.Range("E30").Value = "  +0.97%  "

$ws.Range("E31").Value = "  -3.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.10"
$ws.Range("E32").Value = "  -2.76%  "

$ws.Range("E33").Value = "  -0.80%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  +0.28%  "

$ws.Range("E36").Value = "  -4.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.931"
$ws.Range("E37").Value = "  +1.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.97"
$ws.Range("E38").Value = "  -1.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.19"
$ws.Range("E39").Value = "  +0.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.57"
$ws.Range("E40").Value = "  -2.61%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.79"
$ws.Range("E41").Value = "  +8.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "149.41"
$ws.Range("E42").Value = "  -0.80%  "

$ws.Range("E43").Value = "  -3.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.61"
$ws.Range("E44").Value = "  -1.72%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "283.56"
$ws.Range("E45").Value = "  -2.28%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0929"
$ws.Range("E46").Value = "  -1.24%  "

$ws.Range("E47").Value = "  -1.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.72"
$ws.Range("E48").Value = "  +2.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.559"
$ws.Range("E49").Value = "  -1.78%  "

$ws.Range("E50").Value = "  -1.35%  "

$ws.Range("E51").Value = "  -1.60%  "
